$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00158971523361902
$ws.Range("C2").Value = [double]"6.91180536356096e-05"
$ws.Range("G2").Value = 0.999792645839093
$ws.Range("H2").Value = 0.442977605750622
$ws.Range("I2").Value = 0.999308819463644
$ws.Range("J2").Value = 0.000207354160906829
$ws.Range("L2").Value = 0.996474979264584
$ws.Range("M2").Value = [double]"6.91180536356096e-05"
$ws.Range("N2").Value = 0.000552944429084877
$ws.Range("O2").Value = 0.000276472214542438
$ws.Range("P2").Value = 0.000207354160906829
$ws.Range("Q2").Value = 0.000483826375449267
$ws.Range("R2").Value = 0.000207354160906829
$ws.Range("S2").Value = 0.000898534697262925
$ws.Range("T2").Value = 0.000345590268178048
$ws.Range("U2").Value = 0.000276472214542438
$ws.Range("V2").Value = 0.996405861210948
$ws.Range("W2").Value = 0.915883328725463
$ws.Range("X2").Value = 0.000967652750898535
$ws.Range("B3").Value = 0.997649986176389
$ws.Range("C3").Value = 0.999723527785458
$ws.Range("D3").Value = 0.999239701410008
$ws.Range("E3").Value = 0.000276472214542438
$ws.Range("F3").Value = 0.999792645839093
$ws.Range("H3").Value = [double]"6.91180536356096e-05"
$ws.Range("I3").Value = 0.000276472214542438
$ws.Range("J3").Value = 0.000483826375449267
$ws.Range("K3").Value = 0.99149847940282
$ws.Range("L3").Value = 0.00138236107271219
$ws.Range("M3").Value = 0.998755875034559
$ws.Range("N3").Value = 0.000276472214542438
$ws.Range("O3").Value = 0.998963229195466
$ws.Range("P3").Value = 0.000760298589991706
$ws.Range("Q3").Value = 0.000829416643627315
$ws.Range("R3").Value = 0.000552944429084877
$ws.Range("T3").Value = 0.000138236107271219
$ws.Range("U3").Value = 0.999032347249102
$ws.Range("V3").Value = 0.000483826375449267
$ws.Range("W3").Value = 0.000138236107271219
$ws.Range("X3").Value = 0.998686756980923
$ws.Range("B4").Value = [double]"6.91180536356096e-05"
$ws.Range("F4").Value = [double]"6.91180536356096e-05"
$ws.Range("G4").Value = [double]"6.91180536356096e-05"
$ws.Range("H4").Value = 0.548106165330384
$ws.Range("I4").Value = 0.000414708321813658
$ws.Range("J4").Value = 0.000345590268178048
$ws.Range("L4").Value = 0.000967652750898535
$ws.Range("M4").Value = 0.000552944429084877
$ws.Range("O4").Value = [double]"6.91180536356096e-05"
$ws.Range("Q4").Value = [double]"6.91180536356096e-05"
$ws.Range("R4").Value = 0.000483826375449267
$ws.Range("S4").Value = 0.999032347249102
$ws.Range("T4").Value = 0.999447055570915
$ws.Range("U4").Value = 0.000345590268178048
$ws.Range("V4").Value = 0.00304119435996682
$ws.Range("W4").Value = 0.080729886646392
$ws.Range("X4").Value = 0.000138236107271219
$ws.Range("B5").Value = 0.000414708321813658
$ws.Range("C5").Value = 0.000138236107271219
$ws.Range("D5").Value = 0.000552944429084877
$ws.Range("E5").Value = 0.999516173624551
$ws.Range("F5").Value = 0.000138236107271219
$ws.Range("G5").Value = [double]"6.91180536356096e-05"
$ws.Range("H5").Value = 0.000138236107271219
$ws.Range("J5").Value = 0.998755875034559
$ws.Range("K5").Value = 0.00836328448990876
$ws.Range("L5").Value = 0.000829416643627315
$ws.Range("M5").Value = 0.000207354160906829
$ws.Range("N5").Value = 0.999101465302737
$ws.Range("O5").Value = 0.000207354160906829
$ws.Range("P5").Value = 0.998963229195466
$ws.Range("Q5").Value = 0.998548520873652
$ws.Range("R5").Value = 0.998548520873652
$ws.Range("S5").Value = [double]"6.91180536356096e-05"
$ws.Range("T5").Value = [double]"6.91180536356096e-05"
$ws.Range("U5").Value = 0.000207354160906829
$ws.Range("W5").Value = [double]"6.91180536356096e-05"
